$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Bump the version label on the title page: "V 1.5" -> "V 1.6"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("V 1.5  ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "V 1.6  ", 2) | Out-Null

# ---------------------------------------------------------------------------
# Locate the "Propuesta de aplicación:" section heading. It appears once as
# a short (real) paragraph; the longer hits earlier in the doc belong to the
# table of contents field and are skipped via the length check.
# ---------------------------------------------------------------------------
$headingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*Propuesta de aplicaci*" -and $t.Length -lt 50) {
        $headingIndex = $i
        break
    }
}

if ($headingIndex -eq -1) {
    throw "Could not locate the 'Propuesta de aplicación:' heading paragraph"
}

# The heading is followed by one blank paragraph, then the two paragraphs we
# need to touch:
#   blankIndex + 1 -> tab + "Desarrollaremos un prototipo..." sentence
#   blankIndex + 2 -> four tab characters
$paraDesarrollo = $d.Paragraphs.Item($headingIndex + 2)
$paraTabs       = $d.Paragraphs.Item($headingIndex + 3)

# ---------------------------------------------------------------------------
# 2) Remove the old "Desarrollaremos un prototipo..." sentence, leaving the
#    paragraph's leading tab character untouched.
# ---------------------------------------------------------------------------
$oldSentence = "Desarrollaremos un prototipo que cumpla las funcionalidades que nos ha pedido el cliente. A través de este prototipo el cliente tendrá una idea más clara de cómo va a ser su aplicación y puede que quiera realizar cambios. Del desarrollo completo de la aplicación se encargará otro equipo."
$r1 = $paraDesarrollo.Range
$found1 = $r1.Find.Execute($oldSentence, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 2)
if (-not $found1) {
    throw "Could not find the 'Desarrollaremos un prototipo...' sentence"
}

# ---------------------------------------------------------------------------
# 3) Replace the last three of the four tab characters in the following
#    paragraph with the new, longer description, keeping the first tab.
# ---------------------------------------------------------------------------
$newSentence = "Para que el cliente tenga una idea de la aplicación desarrollaremos la publicación, búsqueda, cancelación y verificación de pedidos y las funcionalidades de valoración entre usuarios y los filtros de búsqueda por valoraciones y cercanía entre ambos, además de bloquear y denunciar para que el cliente tenga una idea de cómo será la aplicación completa. Incluiremos en la implementación también la funcionalidad de usuario premium y usuario básico, al que se le mostrarán pequeños anuncios en la aplicación. También se implementará las publicaciones en twitter a través de la aplicación. Ya que la idea es bastante ambiciosa, símplemente desarrollaremos una pequeña demo de la aplicación para que el cliente tenga una idea de cómo será la aplicación definitiva. En caso de que sobrara tiempo, se desarrollaría el chat y el sistema de log-in y registro."
$tabsRange = $paraTabs.Range
$replaceStart = $tabsRange.Start + 1
$replaceEnd = $tabsRange.End - 1
$r2 = $d.Range($replaceStart, $replaceEnd)
$r2.Text = $newSentence
